$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing keyword in row 37 (was "~", now "(XXL~JL)")
$ws.Range("A37").Value = "(XXL~JL)"

# Append new keyword rows 41-46
$ws.Range("A41").Value = "(25겨울)"
$ws.Range("A42").Value = "(13~19)"
$ws.Range("A43").Value = "(XXL~3XL)"
$ws.Range("A44").Value = "(S~XXL)"
$ws.Range("A45").Value = "(7~11)"
$ws.Range("A46").Value = "(17~19)"
